$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 4900
$ws.Range("I69").Value = 6500
$ws.Range("J69").Value = 3833.3333
$ws.Range("K69").Value = 19500
$ws.Range("L69").Value = 11499.9999
$ws.Range("M69").Value = -18626
$ws.Range("N69").Value = -13247.9999

$ws.Range("H72").Value = 4900
$ws.Range("I72").Value = 6500
$ws.Range("J72").Value = 3833.3333
$ws.Range("K72").Value = 58500
$ws.Range("L72").Value = 34499.9997
$ws.Range("M72").Value = -54132
$ws.Range("N72").Value = -43235.9997

$ws.Range("H103").Value = 7513038
$ws.Range("J103").Value = 740
$ws.Range("L103").Value = 2220
$ws.Range("N103").Value = -3392

$ws.Range("H137").Value = 7143838.5
$ws.Range("I137").Value = 986.6667
$ws.Range("J137").Value = 50000948
$ws.Range("K137").Value = 2960.0001
$ws.Range("L137").Value = 150002844
$ws.Range("M137").Value = -410.0001000000002
$ws.Range("N137").Value = -150007944

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 16668616
$ws.Range("I61").Value = 18520428
$ws.Range("K61").Value = 18520428
$ws.Range("M61").Value = -18520216

$ws.Range("H74").Value = 11907555
$ws.Range("I74").Value = 17243054
$ws.Range("J74").Value = 5288.846
$ws.Range("K74").Value = 17243054
$ws.Range("L74").Value = 5288.846
$ws.Range("M74").Value = -17242180
$ws.Range("N74").Value = -7036.846

$ws.Range("H77").Value = 11907555
$ws.Range("I77").Value = 17243054
$ws.Range("J77").Value = 5288.846
$ws.Range("K77").Value = 86215270
$ws.Range("L77").Value = 26444.23
$ws.Range("M77").Value = -86210902
$ws.Range("N77").Value = -35180.23

$ws.Range("H132").Value = 19235100
$ws.Range("I132").Value = 25003630
$ws.Range("J132").Value = 6670.6665
$ws.Range("K132").Value = 75010890
$ws.Range("L132").Value = 20011.9995
$ws.Range("M132").Value = -75008360
$ws.Range("N132").Value = -25071.9995

$ws.Range("H136").Value = 16668616
$ws.Range("I136").Value = 18520428
$ws.Range("K136").Value = 55561284
$ws.Range("M136").Value = -55558734

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3152.4707
$ws.Range("I134").Value = 2109.7932
$ws.Range("J134").Value = 9200
$ws.Range("K134").Value = 6329.3796
$ws.Range("L134").Value = 27600
$ws.Range("M134").Value = -3794.3796
$ws.Range("N134").Value = -32670

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 15153037
$ws.Range("I31").Value = 1593.8572
$ws.Range("K31").Value = 1593.8572
$ws.Range("M31").Value = -1298.8572

$ws.Range("H34").Value = 15153037
$ws.Range("I34").Value = 1593.8572
$ws.Range("K34").Value = 1593.8572
$ws.Range("M34").Value = -1391.8572

$ws.Range("H57").Value = 22000
$ws.Range("J57").Value = 22000
$ws.Range("L57").Value = 22000
$ws.Range("N57").Value = -23120

$ws.Range("H58").Value = 1343.902
$ws.Range("I58").Value = 616.51514
$ws.Range("J58").Value = 2677.4443
$ws.Range("K58").Value = 616.51514
$ws.Range("L58").Value = 2677.4443
$ws.Range("M58").Value = -413.51514
$ws.Range("N58").Value = -3083.4443

$ws.Range("H132").Value = 38466290
$ws.Range("I132").Value = 83337810
$ws.Range("J132").Value = 4989.4287
$ws.Range("K132").Value = 250013430
$ws.Range("L132").Value = 14968.2861
$ws.Range("M132").Value = -250010900
$ws.Range("N132").Value = -20028.2861

$ws.Range("H134").Value = 1192200.8
$ws.Range("I134").Value = 1860.6111
$ws.Range("J134").Value = 11905262
$ws.Range("K134").Value = 5581.8333
$ws.Range("L134").Value = 35715786
$ws.Range("M134").Value = -3046.8333
$ws.Range("N134").Value = -35720856

$ws.Range("H136").Value = 1343.902
$ws.Range("I136").Value = 616.51514
$ws.Range("J136").Value = 2677.4443
$ws.Range("K136").Value = 1849.54542
$ws.Range("L136").Value = 8032.3329
$ws.Range("M136").Value = 700.4545800000001
$ws.Range("N136").Value = -13132.3329

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 503.86206
$ws.Range("J39").Value = 503.86206
$ws.Range("L39").Value = 1511.58618
$ws.Range("N39").Value = -2099.58618

$ws.Range("H86").Value = 2342.8572
$ws.Range("I86").Value = 1000
$ws.Range("J86").Value = 2880
$ws.Range("K86").Value = 3000
$ws.Range("L86").Value = 8640
$ws.Range("M86").Value = -1814
$ws.Range("N86").Value = -11012

$ws.Range("H89").Value = 2342.8572
$ws.Range("I89").Value = 1000
$ws.Range("J89").Value = 2880
$ws.Range("K89").Value = 9000
$ws.Range("L89").Value = 25920
$ws.Range("M89").Value = -3072
$ws.Range("N89").Value = -37776

$ws.Range("H131").Value = 809.1900000000001
$ws.Range("I131").Value = 395
$ws.Range("J131").Value = 855.2111
$ws.Range("K131").Value = 1185
$ws.Range("L131").Value = 2565.6333
$ws.Range("M131").Value = 3855
$ws.Range("N131").Value = -12645.6333

$ws.Range("H134").Value = 3648.5715
$ws.Range("I134").Value = 1770.7693
$ws.Range("J134").Value = 6700
$ws.Range("K134").Value = 5312.3079
$ws.Range("L134").Value = 20100
$ws.Range("M134").Value = -242.3078999999998
$ws.Range("N134").Value = -30240

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3214.82
$ws.Range("I132").Value = 1992.5
$ws.Range("J132").Value = 7085.5
$ws.Range("K132").Value = 5977.5
$ws.Range("L132").Value = 21256.5
$ws.Range("M132").Value = -3447.5
$ws.Range("N132").Value = -26316.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 6583.778
$ws.Range("I122").Value = 6957.091
$ws.Range("J122").Value = 5997.143
$ws.Range("K122").Value = 20871.273
$ws.Range("L122").Value = 17991.429
$ws.Range("M122").Value = -18421.273
$ws.Range("N122").Value = -22891.429

$ws.Range("H132").Value = 9440660
$ws.Range("I132").Value = 4039.2856
$ws.Range("J132").Value = 20009674
$ws.Range("K132").Value = 12117.8568
$ws.Range("L132").Value = 60029022
$ws.Range("M132").Value = -9587.856800000001
$ws.Range("N132").Value = -60034082

$ws.Range("H136").Value = 15631013
$ws.Range("I136").Value = 29414470
$ws.Range("J136").Value = 9761.4
$ws.Range("K136").Value = 88243410
$ws.Range("L136").Value = 29284.2
$ws.Range("M136").Value = -88240860
$ws.Range("N136").Value = -34384.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2259.7307
$ws.Range("I132").Value = 1491.875
$ws.Range("J132").Value = 3488.3
$ws.Range("K132").Value = 4475.625
$ws.Range("L132").Value = 10464.9
$ws.Range("M132").Value = -1945.625
$ws.Range("N132").Value = -15524.9

$ws.Range("H136").Value = 1167.0476
$ws.Range("I136").Value = 1200.4
$ws.Range("J136").Value = 500
$ws.Range("K136").Value = 3601.2
$ws.Range("L136").Value = 1500
$ws.Range("M136").Value = -1051.2
$ws.Range("N136").Value = -6600
